$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Shapes.Item(2).TextFrame.TextRange.Text = "20 Maggio, 2024"
